$wb = $excel.ActiveWorkbook
$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

$ws_ALC.Range("H62").Value = 2670
$ws_ALC.Range("I62").Value = 2643.5833
$ws_ALC.Range("J62").Value = 2749.25
$ws_ALC.Range("K62").Value = 2643.5833
$ws_ALC.Range("L62").Value = 2749.25
$ws_ALC.Range("M62").Value = -2019.5833
$ws_ALC.Range("N62").Value = -3997.25
$ws_ALC.Range("H65").Value = 2670
$ws_ALC.Range("I65").Value = 2643.5833
$ws_ALC.Range("J65").Value = 2749.25
$ws_ALC.Range("K65").Value = 13217.9165
$ws_ALC.Range("L65").Value = 13746.25
$ws_ALC.Range("M65").Value = -10097.9165
$ws_ALC.Range("N65").Value = -19986.25
$ws_ALC.Range("H74").Value = 11845
$ws_ALC.Range("I74").Value = 15267.5
$ws_ALC.Range("J74").Value = 5000
$ws_ALC.Range("K74").Value = 15267.5
$ws_ALC.Range("L74").Value = 5000
$ws_ALC.Range("M74").Value = -14331.5
$ws_ALC.Range("N74").Value = -6872
$ws_ALC.Range("H77").Value = 11845
$ws_ALC.Range("I77").Value = 15267.5
$ws_ALC.Range("J77").Value = 5000
$ws_ALC.Range("K77").Value = 76337.5
$ws_ALC.Range("L77").Value = 25000
$ws_ALC.Range("M77").Value = -71657.5
$ws_ALC.Range("N77").Value = -34360
$ws_ALC.Range("H133").Value = 37666.668
$ws_ALC.Range("J133").Value = 37666.668
$ws_ALC.Range("L133").Value = 37666.668
$ws_ALC.Range("N133").Value = -47786.668
$ws_ALC.Range("H134").Value = 34780
$ws_ALC.Range("J134").Value = 34780
$ws_ALC.Range("L134").Value = 34780
$ws_ALC.Range("N134").Value = -44920
$ws_ALC.Range("H136").Value = 59933
$ws_ALC.Range("J136").Value = 59933
$ws_ALC.Range("L136").Value = 59933
$ws_ALC.Range("N136").Value = -70133
$ws_ALC.Range("H138").Value = 3921.7273
$ws_ALC.Range("I138").Value = 2585.7334
$ws_ALC.Range("J138").Value = 4612.759
$ws_ALC.Range("K138").Value = 7757.2002
$ws_ALC.Range("L138").Value = 13838.277
$ws_ALC.Range("M138").Value = -2617.2002
$ws_ALC.Range("N138").Value = -24118.277
$ws_ALC.Range("H139").Value = 59999.75
$ws_ALC.Range("J139").Value = 69999.664
$ws_ALC.Range("L139").Value = 69999.664
$ws_ALC.Range("N139").Value = -80279.664
$ws_ALC.Range("H140").Value = 0
$ws_ALC.Range("J140").Value = 0
$ws_ALC.Range("L140").Value = 0
$ws_ALC.Range("N140").ClearContents()
$ws_ARM.Range("H32").Value = 30482.8
$ws_ARM.Range("I32").Value = 11061.358
$ws_ARM.Range("J32").Value = 130577.92
$ws_ARM.Range("K32").Value = 11061.358
$ws_ARM.Range("L32").Value = 130577.92
$ws_ARM.Range("M32").Value = -10774.358
$ws_ARM.Range("N32").Value = -131151.92
$ws_ARM.Range("H123").Value = 50000
$ws_ARM.Range("J123").Value = 50000
$ws_ARM.Range("L123").Value = 50000
$ws_ARM.Range("N123").Value = -59800
$ws_ARM.Range("H135").Value = 45753.168
$ws_ARM.Range("J135").Value = 45753.168
$ws_ARM.Range("L135").Value = 45753.168
$ws_ARM.Range("N135").Value = -55893.168
$ws_ARM.Range("H138").Value = 75500
$ws_ARM.Range("J138").Value = 75500
$ws_ARM.Range("L138").Value = 75500
$ws_ARM.Range("N138").Value = -85780
$ws_ARM.Range("H139").Value = 62530.332
$ws_ARM.Range("I139").Value = 0
$ws_ARM.Range("J139").Value = 62530.332
$ws_ARM.Range("K139").Value = 0
$ws_ARM.Range("L139").Value = 62530.332
$ws_ARM.Range("M139").ClearContents()
$ws_ARM.Range("N139").Value = -72810.33199999999
$ws_BSM.Range("H135").Value = 58450
$ws_BSM.Range("J135").Value = 58450
$ws_BSM.Range("L135").Value = 58450
$ws_BSM.Range("N135").Value = -68590
$ws_BSM.Range("H137").Value = 38769.5
$ws_BSM.Range("J137").Value = 38769.5
$ws_BSM.Range("L137").Value = 38769.5
$ws_BSM.Range("N137").Value = -48969.5
$ws_BSM.Range("H138").Value = 134000
$ws_BSM.Range("J138").Value = 134000
$ws_BSM.Range("L138").Value = 134000
$ws_BSM.Range("N138").Value = -144280
$ws_BSM.Range("H140").Value = 49472.5
$ws_BSM.Range("J140").Value = 49472.5
$ws_BSM.Range("L140").Value = 49472.5
$ws_BSM.Range("N140").Value = -59832.5
$ws_CRP.Range("H31").Value = 53796.25
$ws_CRP.Range("I31").Value = 1339.8462
$ws_CRP.Range("K31").Value = 1339.8462
$ws_CRP.Range("M31").Value = -1044.8462
$ws_CRP.Range("H34").Value = 53796.25
$ws_CRP.Range("I34").Value = 1339.8462
$ws_CRP.Range("K34").Value = 1339.8462
$ws_CRP.Range("M34").Value = -1137.8462
$ws_CRP.Range("H138").Value = 149000
$ws_CRP.Range("J138").Value = 149000
$ws_CRP.Range("L138").Value = 149000
$ws_CRP.Range("N138").Value = -159280
$ws_CRP.Range("H140").Value = 49000
$ws_CRP.Range("I140").Value = 30000
$ws_CRP.Range("J140").Value = 68000
$ws_CRP.Range("K140").Value = 30000
$ws_CRP.Range("L140").Value = 68000
$ws_CRP.Range("M140").Value = -24820
$ws_CRP.Range("N140").Value = -78360
$ws_CUL.Range("H56").Value = 4268.5454
$ws_CUL.Range("I56").Value = 4268.5454
$ws_CUL.Range("K56").Value = 4268.5454
$ws_CUL.Range("M56").Value = -3738.5454
$ws_CUL.Range("H113").Value = 1009.53845
$ws_CUL.Range("J113").Value = 647.05884
$ws_CUL.Range("L113").Value = 1941.17652
$ws_CUL.Range("N113").Value = -6281.17652
$ws_GSM.Range("H70").Value = 132012.81
$ws_GSM.Range("I70").Value = 173234
$ws_GSM.Range("J70").Value = 8349.25
$ws_GSM.Range("K70").Value = 173234
$ws_GSM.Range("L70").Value = 8349.25
$ws_GSM.Range("M70").Value = -172964
$ws_GSM.Range("N70").Value = -8889.25
$ws_GSM.Range("H73").Value = 132012.81
$ws_GSM.Range("I73").Value = 173234
$ws_GSM.Range("J73").Value = 8349.25
$ws_GSM.Range("K73").Value = 173234
$ws_GSM.Range("L73").Value = 8349.25
$ws_GSM.Range("M73").Value = -172298
$ws_GSM.Range("N73").Value = -10221.25
$ws_GSM.Range("H102").Value = 2045.5238
$ws_GSM.Range("I102").Value = 1635.375
$ws_GSM.Range("J102").Value = 3358
$ws_GSM.Range("K102").Value = 1635.375
$ws_GSM.Range("L102").Value = 3358
$ws_GSM.Range("M102").Value = -13.375
$ws_GSM.Range("N102").Value = -6602
$ws_GSM.Range("H135").Value = 44626.285
$ws_GSM.Range("J135").Value = 44626.285
$ws_GSM.Range("L135").Value = 44626.285
$ws_GSM.Range("N135").Value = -54766.285
$ws_GSM.Range("H138").Value = 97600
$ws_GSM.Range("I138").Value = 88000
$ws_GSM.Range("J138").Value = 100000
$ws_GSM.Range("K138").Value = 88000
$ws_GSM.Range("L138").Value = 100000
$ws_GSM.Range("M138").Value = -82860
$ws_GSM.Range("N138").Value = -110280
$ws_GSM.Range("H140").Value = 134120
$ws_GSM.Range("J140").Value = 134120
$ws_GSM.Range("L140").Value = 134120
$ws_GSM.Range("N140").Value = -144480
$ws_GSM.Range("H141").Value = 42890
$ws_GSM.Range("J141").Value = 42890
$ws_GSM.Range("L141").Value = 42890
$ws_GSM.Range("N141").Value = -53250
$ws_LTW.Range("H40").Value = 58518.777
$ws_LTW.Range("I40").Value = 169483
$ws_LTW.Range("J40").Value = 3036.6667
$ws_LTW.Range("K40").Value = 169483
$ws_LTW.Range("L40").Value = 3036.6667
$ws_LTW.Range("M40").Value = -169347
$ws_LTW.Range("N40").Value = -3308.6667
$ws_LTW.Range("H136").Value = 2385.9412
$ws_LTW.Range("I136").Value = 1593.2
$ws_LTW.Range("J136").Value = 3518.4285
$ws_LTW.Range("K136").Value = 4779.6
$ws_LTW.Range("L136").Value = 10555.2855
$ws_LTW.Range("M136").Value = -2229.6
$ws_LTW.Range("N136").Value = -15655.2855
$ws_LTW.Range("H139").Value = 54996
$ws_LTW.Range("I139").Value = 16000
$ws_LTW.Range("J139").Value = 64745
$ws_LTW.Range("K139").Value = 16000
$ws_LTW.Range("L139").Value = 64745
$ws_LTW.Range("M139").Value = -10860
$ws_LTW.Range("N139").Value = -75025
$ws_WVR.Range("H58").Value = 13490
$ws_WVR.Range("I58").Value = 3000
$ws_WVR.Range("J58").Value = 23980
$ws_WVR.Range("K58").Value = 3000
$ws_WVR.Range("L58").Value = 23980
$ws_WVR.Range("M58").Value = -2692
$ws_WVR.Range("N58").Value = -24596
$ws_WVR.Range("H132").Value = 13548
$ws_WVR.Range("I132").Value = 11771.786
$ws_WVR.Range("J132").Value = 17100.428
$ws_WVR.Range("K132").Value = 35315.358
$ws_WVR.Range("L132").Value = 51301.284
$ws_WVR.Range("M132").Value = -32785.358
$ws_WVR.Range("N132").Value = -56361.284
$ws_WVR.Range("H137").Value = 41285.832
$ws_WVR.Range("J137").Value = 41285.832
$ws_WVR.Range("L137").Value = 41285.832
$ws_WVR.Range("N137").Value = -51485.832
$ws_WVR.Range("H138").Value = 49431.527
$ws_WVR.Range("J138").Value = 49431.527
$ws_WVR.Range("L138").Value = 49431.527
$ws_WVR.Range("N138").Value = -59711.527
$ws_WVR.Range("H139").Value = 64135
$ws_WVR.Range("J139").Value = 64135
$ws_WVR.Range("L139").Value = 64135
$ws_WVR.Range("N139").Value = -74415
